$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (column B 52->35, column D 41->30) ---
# ColumnWidth setter stores value+5/6 in the file, so subtract 5/6 to land on the exact target.
$ws.Columns(2).ColumnWidth = 35 - 5/6
$ws.Columns(4).ColumnWidth = 30 - 5/6

# --- Row 2: only the captured-at timestamp changes ---
$ws.Range("A2").Value = "2025-12-06 06:25:54"

# --- Row 3: new listing content ---
$ws.Range("A3").Value = "2025-12-06 06:25:54"
$ws.Range("B3").Value = "【急募】新規システム開発に伴う要件定義依頼"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5448563"
$ws.Range("G3").Value = 110
$ws.Range("H3").Value = "◆開発,システム開発"

# --- Row 4: new listing content ---
$ws.Range("A4").Value = "2025-12-06 06:25:54"
$ws.Range("B4").Value = "【受注メールを元にECサイト自動仕入ツール】"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5448428"
$ws.Range("G4").Value = 98
$ws.Range("H4").Value = "◆ツール ◇サイト"

# --- Row 5: new listing content ---
$ws.Range("A5").Value = "2025-12-06 06:25:54"
$ws.Range("B5").Value = "【緊急】既存コードの構造解析ができるPHPエンジニアを探しています"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5448440"
$ws.Range("G5").Value = 33
$ws.Range("H5").Value = "○PHP"

# --- Drop the old rows 6-11 entirely (dimension shrinks to A1:H5) ---
$ws.Rows("6:11").Delete()

# The engine's Hyperlinks collection doesn't auto-clean stale entries (old F6:F11
# refs) when rows are deleted, and single-item .Delete() on a Hyperlink object is a
# no-op here, so clear the whole collection and rebuild only the four that remain.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5448409")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5448563")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5448428")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5448440")
# Re-adding a hyperlink nudges the cell onto a freshly duplicated style; restore the
# original shared "Hyperlink" cell style so column F keeps its original formatting.
$ws.Range("F2:F5").Style = "Hyperlink"
